# Updated symbol list on Fri Jan 13 08:56:27 UTC 2023 with GitHub Actions
# Refresh the Price (D) and Volume(1h) (E) columns for the cryptos sheet.
# Values are entered with a leading apostrophe so Excel keeps them as
# literal text (preserving trailing zeros / the "%" suffix) instead of
# auto-converting them into numeric / percentage cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'287.20"
$ws.Range("E2").Value = "'1.45%"
$ws.Range("D3").Value = "'29.55"
$ws.Range("E3").Value = "'3.41%"
$ws.Range("D4").Value = "'5.114"
$ws.Range("E4").Value = "'1.00%"
$ws.Range("D5").Value = "'0.06698"
$ws.Range("E5").Value = "'3.27%"
$ws.Range("D6").Value = "'7.341"
$ws.Range("E6").Value = "'1.68%"
$ws.Range("E7").Value = "'1.16%"
$ws.Range("D8").Value = "'1.359"
$ws.Range("E8").Value = "'-3.94%"
$ws.Range("D9").Value = "'0.9120"
$ws.Range("E9").Value = "'0.13%"
$ws.Range("D10").Value = "'0.1590"
$ws.Range("E10").Value = "'2.76%"
$ws.Range("D11").Value = "'0.06739"
$ws.Range("E11").Value = "'3.08%"
$ws.Range("D12").Value = "'0.07670"
$ws.Range("E12").Value = "'1.29%"
$ws.Range("D13").Value = "'0.02929"
$ws.Range("E13").Value = "'6.20%"
$ws.Range("D14").Value = "'0.08978"
$ws.Range("E14").Value = "'0.23%"
$ws.Range("D15").Value = "'0.001563"
$ws.Range("E15").Value = "'-1.53%"
$ws.Range("D16").Value = "'0.04485"
$ws.Range("E16").Value = "'0.84%"
$ws.Range("D17").Value = "'0.0006467"
$ws.Range("E17").Value = "'1.27%"
$ws.Range("D18").Value = "'0.006263"
$ws.Range("E18").Value = "'2.05%"
$ws.Range("D19").Value = "'3.442"
$ws.Range("E19").Value = "'-0.24%"
$ws.Range("D20").Value = "'2.220"
$ws.Range("E20").Value = "'-0.97%"
$ws.Range("E21").Value = "'0.76%"
$ws.Range("E22").Value = "'-2.38%"
$ws.Range("D23").Value = "'4.067"
$ws.Range("E23").Value = "'2.33%"
$ws.Range("E24").Value = "'2.42%"
$ws.Range("E25").Value = "'0.46%"
$ws.Range("D26").Value = "'0.004118"
$ws.Range("E26").Value = "'-7.70%"
$ws.Range("E27").Value = "'-0.12%"
$ws.Range("D28").Value = "'0.0001616"
$ws.Range("E28").Value = "'-1.18%"
$ws.Range("D40").Value = "'0.04257"
$ws.Range("E40").Value = "'3.67%"
$ws.Range("D41").Value = "'0.006776"
$ws.Range("E41").Value = "'2.31%"
$ws.Range("E42").Value = "'0.76%"
$ws.Range("D43").Value = "'0.002188"
$ws.Range("E43").Value = "'6.70%"
$ws.Range("D44").Value = "'0.01328"
$ws.Range("E44").Value = "'7.24%"
$ws.Range("D45").Value = "'0.00005674"
$ws.Range("E45").Value = "'4.96%"
$ws.Range("D46").Value = "'1.974"
$ws.Range("E46").Value = "'2.14%"
$ws.Range("E47").Value = "'-29.46%"
